# Update database and change read_price algorithm:
# shift quarterly columns (D..M) one quarter to the left, dropping the
# oldest quarter (Q2 1399/06) and appending the newest quarter
# (Q4 1401/12) with placeholder values, mirroring the upstream data
# refresh captured in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: quarter-period header labels (D8:M8)
$periods = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $periods[$i]
}

# Row 9: publish-date header labels (D9:M9)
$dates = @(
    "1400-11-05 (4)",
    "1401-04-18 (10)",
    "1401-04-21 (2)",
    "1401-08-30 (4)",
    "1401-10-29 (3)",
    "1402-02-27 (7)",
    "1401-04-21",
    "1401-08-30 (2)",
    "1401-10-29",
    "1402-02-27"
)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $dates[$i]
}

# Rows 11-27: financial data rows. The whole D:M block is reset (shifted
# left and cleared for the new quarter) -- every numeric cell becomes 0,
# except row 15 and row 23 ("not applicable" rows) which become the
# text placeholder "-" across the whole row, and row 16 column J which
# also becomes "-".
for ($row = 11; $row -le 27; $row++) {
    for ($col = 4; $col -le 13; $col++) {
        if ($row -eq 15 -or $row -eq 23) {
            $ws.Cells.Item($row, $col).Value = "-"
        } elseif ($row -eq 16 -and $col -eq 10) {
            $ws.Cells.Item($row, $col).Value = "-"
        } else {
            $ws.Cells.Item($row, $col).Value = 0
        }
    }
}
